# Edit script: add a new "数组" (Array) sheet row for LeetCode #219
# "Contains Duplicate II" problem, matching the commit "repeat 2 code with array".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")
$ws.Activate()

$problemText = @'
给定一个整数数组和一个整数 k，判断数组中是否存在两个不同的索引 i 和 j，使得 nums [i] = nums [j]，并且 i 和 j 的差的 绝对值
 至多为 k。 
 示例 1:
 输入: nums = [1,2,3,1], k = 3
输出: true 
 示例 2:
 输入: nums = [1,0,1,1], k = 1
输出: true 
 示例 3:
 输入: nums = [1,2,3,1,2,3], k = 2
输出: false 
'@

$approachText = @'
1 使用hash表存储数组元素以及相应的索引
2 如果不存在，舅加入map
2 如果元素已存在，比较两个索引是否满足条件，满足就返回，不满足需要更新元素在map中的值，处理(1,0,1,1 1)的情况
3 迭代结束条件是元素是否已迭代完成
'@

# New row 10 data
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 219
$ws.Range("C10").Value = $problemText
$ws.Range("D10").Value = $approachText
$ws.Range("E10").Value = "更新map"
$ws.Range("F10").Value = "O(N)"
$ws.Range("G10").Value = "O(N)"

# Row height for the newly added row
$ws.Rows.Item(10).RowHeight = 286

# Update the view: scroll so the new row is visible, select D15 to match the
# post-edit cursor position recorded in the workbook.
$ws.Range("D15").Select()
